$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), R (Origen), S (Precio $/Kg)
# Row 6 is intentionally left untouched (it matches the diff: no hunk touches row 6).

$targets = @{
    2  = @{ D = 44586; M = 80;  N = 7000; O = 7000; P = 7000; R = "Provincia de Curicó"; S = 3500 }
    3  = @{ D = 44588; M = 160; N = 6500; O = 7000; P = 6750; R = "Provincia de Curicó"; S = 3375 }
    4  = @{ D = 44587; M = 165; N = 6500; O = 7000; P = 6742; R = "Provincia de Linares"; S = 3371 }
    5  = @{ D = 44214; M = 48;  N = 6000; O = 6000; P = 6000; R = "Provincia de Linares"; S = 3000 }
    7  = @{ D = 44585; M = 160; N = 6500; O = 7000; P = 6750; R = "Provincia de Curicó"; S = 3375 }
    8  = @{ D = 44211; M = 45;  N = 6000; O = 6000; P = 6000; R = "Provincia de Curicó"; S = 3000 }
    9  = @{ D = 44592; M = 30;  N = 8000; O = 8000; P = 8000; R = "Provincia de Linares"; S = 4000 }
    10 = @{ D = 44582; M = 150; N = 6000; O = 6500; P = 6233; R = "Provincia de Curicó"; S = 3116 }
    11 = @{ D = 44614; M = 45;  N = 6000; O = 6000; P = 6000; R = "Provincia de Linares"; S = 3000 }
    12 = @{ D = 44589; M = 60;  N = 6000; O = 6000; P = 6000; R = "Provincia de Curicó"; S = 3000 }
    13 = @{ D = 44628; M = 40;  N = 6000; O = 6000; P = 6000; R = "Provincia de Linares"; S = 3000 }
    14 = @{ D = 44209; M = 58;  N = 6000; O = 6000; P = 6000; R = "Provincia de Curicó"; S = 3000 }
    15 = @{ D = 44627; M = 45;  N = 6000; O = 6000; P = 6000; R = "Provincia de Linares"; S = 3000 }
}

foreach ($row in $targets.Keys) {
    $t = $targets[$row]
    $ws.Cells.Item($row, 4).Value = $t.D    # D - Fecha
    $ws.Cells.Item($row, 13).Value = $t.M   # M - Volumen
    $ws.Cells.Item($row, 14).Value = $t.N   # N - Precio minimo
    $ws.Cells.Item($row, 15).Value = $t.O   # O - Precio maximo
    $ws.Cells.Item($row, 16).Value = $t.P   # P - Precio promedio ponderado
    $ws.Cells.Item($row, 18).Value = $t.R   # R - Origen
    $ws.Cells.Item($row, 19).Value = $t.S   # S - Precio $/Kg
}
